# The commit re-runs the "insert each dictionary row into the DB" script one
# more time against the same worksheet. Since the sheet stores the running
# row index in column A, every index shifts up by one (0->1, 1->2, 2->3,
# 3->4, 4->5, 5->6) and the active selection moves on to the next empty
# staging row (E8) ready for the following run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump the running counter in column A by one for every data row.
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 6

# A3 and A5 pick up the same (no-op) "General" alignment flag the header
# row (A1) already carries, matching how the other rows look after the
# re-run.
$ws.Range("A3").HorizontalAlignment = 1
$ws.Range("A5").HorizontalAlignment = 1

# Leave the selection parked on the next input row for the following
# dictionary entry.
$ws.Range("E8").Select()
